# Auto-generated edit script: updates crypto price/volume table
# Source data (coinranking.com) refreshed; several coins also swapped rank
# position (rows 24/25 and 44/45) between the two scrapes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    # Writes $Val into (Row, Col) and makes sure it survives as TEXT,
    # matching the source file where every Coin/Link/Price/Volume cell
    # is an inline string (e.g. '1.00' must not collapse to the number 1,
    # '65.753.27' must not be misread as a multi-part number, etc).
    param($Row, $Col, $Val)
    $c = $ws.Cells.Item($Row, $Col)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = $origStyle
}

Set-CellText 2 4 '65.753.27'
Set-CellText 2 5 '  +0.14%  '

Set-CellText 3 4 '2.653.81'
Set-CellText 3 5 '  -0.62%  '

Set-CellText 4 5 '  +0.13%  '

Set-CellText 5 4 '597.53'
Set-CellText 5 5 '  -0.81%  '

Set-CellText 6 4 '157.07'
Set-CellText 6 5 '  -0.29%  '

Set-CellText 7 5 '  +0.12%  '

Set-CellText 8 4 '0.630'
Set-CellText 8 5 '  +1.40%  '

Set-CellText 9 4 '0.126'
Set-CellText 9 5 '  +0.98%  '

Set-CellText 10 5 '  -1.13%  '

Set-CellText 11 5 '  -0.89%  '

Set-CellText 12 5 '  +1.17%  '

Set-CellText 13 4 '28.62'
Set-CellText 13 5 '  -2.79%  '

Set-CellText 14 5 '  -0.23%  '

Set-CellText 15 4 '3.131.77'
Set-CellText 15 5 '  -0.32%  '

Set-CellText 16 4 '65.580.36'
Set-CellText 16 5 '  +0.16%  '

Set-CellText 17 4 '2.637.89'
Set-CellText 17 5 '  -0.97%  '

Set-CellText 18 4 '12.55'
Set-CellText 18 5 '  -1.21%  '

Set-CellText 19 5 '  -1.50%  '

Set-CellText 20 4 '7.46'
Set-CellText 20 5 '  -2.69%  '

Set-CellText 21 4 '349.47'
Set-CellText 21 5 '  -0.49%  '

Set-CellText 22 5 '  +0.18%  '

Set-CellText 23 4 '69.24'
Set-CellText 23 5 '  -0.31%  '

Set-CellText 24 2 'SuiNetwork'
Set-CellText 24 3 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-CellText 24 4 '1.76'
Set-CellText 24 5 '  +8.44%  '

Set-CellText 25 2 'PEPE'
Set-CellText 25 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText 25 4 '0.0000112'
Set-CellText 25 5 '  +1.46%  '

Set-CellText 26 4 '9.61'
Set-CellText 26 5 '  -1.33%  '

Set-CellText 27 5 '  +0.76%  '

Set-CellText 28 4 '561.57'
Set-CellText 28 5 '  +5.67%  '

Set-CellText 29 4 '0.163'
Set-CellText 29 5 '  -2.31%  '

Set-CellText 30 4 '1.00'
Set-CellText 30 5 '  -0.04%  '

Set-CellText 31 4 '7.91'
Set-CellText 31 5 '  -1.82%  '

Set-CellText 32 4 '2.13'
Set-CellText 32 5 '  -0.21%  '

Set-CellText 33 5 '  +1.12%  '

Set-CellText 34 5 '  -0.75%  '

Set-CellText 35 5 '  -1.04%  '

Set-CellText 36 5 '  -0.91%  '

Set-CellText 37 4 '20.39'
Set-CellText 37 5 '  -0.03%  '

Set-CellText 38 4 '1.00'
Set-CellText 38 5 '  +0.06%  '

Set-CellText 39 5 '  -0.22%  '

Set-CellText 40 4 '154.92'
Set-CellText 40 5 '  -2.70%  '

Set-CellText 41 5 '  -0.02%  '

Set-CellText 42 4 '161.71'
Set-CellText 42 5 '  -2.01%  '

Set-CellText 43 5 '  -0.39%  '

Set-CellText 44 2 'Hedera'
Set-CellText 44 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText 44 4 '0.0606'
Set-CellText 44 5 '  -0.86%  '

Set-CellText 45 2 'dogwifhat'
Set-CellText 45 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText 45 4 '2.28'
Set-CellText 45 5 '  -0.62%  '

Set-CellText 46 4 '22.62'
Set-CellText 46 5 '  -1.76%  '

Set-CellText 47 4 '0.637'
Set-CellText 47 5 '  -1.03%  '

Set-CellText 48 5 '  -1.46%  '

Set-CellText 49 5 '  -0.95%  '

Set-CellText 50 4 '19.75'
Set-CellText 50 5 '  -2.35%  '

Set-CellText 51 5 '  +6.78%  '
